$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2311.5
$ws.Range("J98").Value = 482.25
$ws.Range("L98").Value = 482.25
$ws.Range("N98").Value = -3478.25

$ws.Range("H122").Value = 2311.5
$ws.Range("J122").Value = 482.25
$ws.Range("L122").Value = 1446.75
$ws.Range("N122").Value = -6346.75

$ws.Range("H132").Value = 388247.88
$ws.Range("I132").Value = 403757.56
$ws.Range("J132").Value = 506
$ws.Range("K132").Value = 1211272.68
$ws.Range("L132").Value = 1518
$ws.Range("M132").Value = -1208742.68
$ws.Range("N132").Value = -6578

$ws.Range("H137").Value = 37039000
$ws.Range("I137").Value = 1289.2778
$ws.Range("K137").Value = 3867.8334
$ws.Range("M137").Value = -1317.8334

$ws.Range("H138").Value = 1796.3914
$ws.Range("I138").Value = 1203.6316
$ws.Range("J138").Value = 2523
$ws.Range("K138").Value = 3610.8948
$ws.Range("L138").Value = 7569
$ws.Range("M138").Value = 1529.1052
$ws.Range("N138").Value = -17849

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 883.1
$ws.Range("I45").Value = 900.3333
$ws.Range("J45").Value = 875.7143
$ws.Range("K45").Value = 900.3333
$ws.Range("L45").Value = 875.7143
$ws.Range("M45").Value = -523.3333
$ws.Range("N45").Value = -1629.7143

$ws.Range("H61").Value = 3325.84
$ws.Range("I61").Value = 2417.8823
$ws.Range("K61").Value = 2417.8823
$ws.Range("M61").Value = -2205.8823

$ws.Range("H74").Value = 3886.6316
$ws.Range("I74").Value = 783.7895
$ws.Range("J74").Value = 6989.4736
$ws.Range("K74").Value = 783.7895
$ws.Range("L74").Value = 6989.4736
$ws.Range("M74").Value = 90.21050000000002
$ws.Range("N74").Value = -8737.473600000001

$ws.Range("H77").Value = 3886.6316
$ws.Range("I77").Value = 783.7895
$ws.Range("J77").Value = 6989.4736
$ws.Range("K77").Value = 3918.9475
$ws.Range("L77").Value = 34947.368
$ws.Range("M77").Value = 449.0525000000002
$ws.Range("N77").Value = -43683.368

$ws.Range("H122").Value = 2999.5
$ws.Range("I122").Value = 2999.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8998.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = -6548.5

$ws.Range("H132").Value = 2270.4194
$ws.Range("I132").Value = 2195.2083
$ws.Range("K132").Value = 6585.624899999999
$ws.Range("M132").Value = -4055.624899999999

$ws.Range("H136").Value = 3325.84
$ws.Range("I136").Value = 2417.8823
$ws.Range("K136").Value = 7253.646900000001
$ws.Range("M136").Value = -4703.646900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1684.2222
$ws.Range("I86").Value = 1488.5938
$ws.Range("K86").Value = 1488.5938
$ws.Range("M86").Value = -365.5938000000001

$ws.Range("H89").Value = 1684.2222
$ws.Range("I89").Value = 1488.5938
$ws.Range("K89").Value = 7442.969000000001
$ws.Range("M89").Value = -1826.969000000001

$ws.Range("H105").Value = 2783.3333
$ws.Range("I105").Value = 3233.3333
$ws.Range("K105").Value = 3233.3333
$ws.Range("M105").Value = -1486.3333

$ws.Range("H107").Value = 2615.0715
$ws.Range("I107").Value = 2118.5
$ws.Range("J107").Value = 2987.5
$ws.Range("K107").Value = 2118.5
$ws.Range("L107").Value = 2987.5
$ws.Range("M107").Value = -198.5
$ws.Range("N107").Value = -6827.5

$ws.Range("H134").Value = 36512.906
$ws.Range("I134").Value = 39685.69
$ws.Range("K134").Value = 119057.07
$ws.Range("M134").Value = -116522.07

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1467.0667
$ws.Range("J31").Value = 1940
$ws.Range("L31").Value = 1940
$ws.Range("N31").Value = -2530

$ws.Range("H34").Value = 1467.0667
$ws.Range("J34").Value = 1940
$ws.Range("L34").Value = 1940
$ws.Range("N34").Value = -2344

$ws.Range("H58").Value = 3607.3157
$ws.Range("I58").Value = 3443.5
$ws.Range("J58").Value = 3682.923
$ws.Range("K58").Value = 3443.5
$ws.Range("L58").Value = 3682.923
$ws.Range("M58").Value = -3240.5
$ws.Range("N58").Value = -4088.923

$ws.Range("H107").Value = 3639.8572
$ws.Range("I107").Value = 4853.2
$ws.Range("J107").Value = 606.5
$ws.Range("K107").Value = 4853.2
$ws.Range("L107").Value = 606.5
$ws.Range("M107").Value = -2933.2
$ws.Range("N107").Value = -4446.5

$ws.Range("H132").Value = 2920.8276
$ws.Range("I132").Value = 2689.3684
$ws.Range("K132").Value = 8068.1052
$ws.Range("M132").Value = -5538.1052

$ws.Range("H134").Value = 2574.625
$ws.Range("I134").Value = 2610.3635
$ws.Range("J134").Value = 2496
$ws.Range("K134").Value = 7831.0905
$ws.Range("L134").Value = 7488
$ws.Range("M134").Value = -5296.0905
$ws.Range("N134").Value = -12558

$ws.Range("H136").Value = 3607.3157
$ws.Range("I136").Value = 3443.5
$ws.Range("J136").Value = 3682.923
$ws.Range("K136").Value = 10330.5
$ws.Range("L136").Value = 11048.769
$ws.Range("M136").Value = -7780.5
$ws.Range("N136").Value = -16148.769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3196.3948
$ws.Range("I113").Value = 530.0476
$ws.Range("J113").Value = 6490.1177
$ws.Range("K113").Value = 1590.1428
$ws.Range("L113").Value = 19470.3531
$ws.Range("M113").Value = 579.8571999999999
$ws.Range("N113").Value = -23810.3531

$ws.Range("H133").Value = 5946.8335
$ws.Range("I133").Value = 3507.1428
$ws.Range("J133").Value = 7499.364
$ws.Range("K133").Value = 10521.4284
$ws.Range("L133").Value = 22498.092
$ws.Range("M133").Value = -5461.428400000001
$ws.Range("N133").Value = -32618.092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 70.625
$ws.Range("J2").Value = 136
$ws.Range("L2").Value = 136
$ws.Range("N2").Value = -362

$ws.Range("H43").Value = 17705.883
$ws.Range("J43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20302

$ws.Range("H46").Value = 16450
$ws.Range("J46").Value = 16450
$ws.Range("L46").Value = 16450
$ws.Range("N46").Value = -16762

$ws.Range("H57").Value = 15074.75
$ws.Range("I57").Value = 10000
$ws.Range("J57").Value = 16766.334
$ws.Range("K57").Value = 10000
$ws.Range("L57").Value = 16766.334
$ws.Range("M57").Value = -9180
$ws.Range("N57").Value = -18406.334

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = ""
$ws.Range("N93").Value = 0

$ws.Range("H122").Value = 3330.0908
$ws.Range("I122").Value = 3878.875
$ws.Range("J122").Value = 1866.6666
$ws.Range("K122").Value = 11636.625
$ws.Range("L122").Value = 5599.9998
$ws.Range("M122").Value = -9186.625
$ws.Range("N122").Value = -10499.9998

$ws.Range("H132").Value = 2526
$ws.Range("I132").Value = 2164.0588
$ws.Range("J132").Value = 2999.3076
$ws.Range("K132").Value = 6492.176399999999
$ws.Range("L132").Value = 8997.9228
$ws.Range("M132").Value = -3962.176399999999
$ws.Range("N132").Value = -14057.9228

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1286.5518
$ws.Range("I61").Value = 1023.55
$ws.Range("J61").Value = 1871
$ws.Range("K61").Value = 1023.55
$ws.Range("L61").Value = 1871
$ws.Range("M61").Value = -821.55
$ws.Range("N61").Value = -2275

$ws.Range("H82").Value = 1295.6538
$ws.Range("I82").Value = 1192.6923
$ws.Range("J82").Value = 1398.6154
$ws.Range("K82").Value = 1192.6923
$ws.Range("L82").Value = 1398.6154
$ws.Range("M82").Value = -831.6922999999999
$ws.Range("N82").Value = -2120.6154

$ws.Range("H85").Value = 1295.6538
$ws.Range("I85").Value = 1192.6923
$ws.Range("J85").Value = 1398.6154
$ws.Range("K85").Value = 1192.6923
$ws.Range("L85").Value = 1398.6154
$ws.Range("M85").Value = 55.30770000000007
$ws.Range("N85").Value = -3894.6154

$ws.Range("H113").Value = 1286.5518
$ws.Range("I113").Value = 1023.55
$ws.Range("J113").Value = 1871
$ws.Range("K113").Value = 1023.55
$ws.Range("L113").Value = 1871
$ws.Range("M113").Value = 1146.45
$ws.Range("N113").Value = -6211

$ws.Range("H122").Value = 5603.4414
$ws.Range("I122").Value = 6547.3335
$ws.Range("J122").Value = 4078.6924
$ws.Range("K122").Value = 19642.0005
$ws.Range("L122").Value = 12236.0772
$ws.Range("M122").Value = -17192.0005
$ws.Range("N122").Value = -17136.0772

$ws.Range("H132").Value = 8300.263000000001
$ws.Range("I132").Value = 9367.134
$ws.Range("J132").Value = 4299.5
$ws.Range("K132").Value = 28101.402
$ws.Range("L132").Value = 12898.5
$ws.Range("M132").Value = -25571.402
$ws.Range("N132").Value = -17958.5

$ws.Range("H136").Value = 2825
$ws.Range("I136").Value = 1266.6666
$ws.Range("J136").Value = 3760
$ws.Range("K136").Value = 3799.9998
$ws.Range("L136").Value = 11280
$ws.Range("M136").Value = -1249.9998
$ws.Range("N136").Value = -16380

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 31983
$ws.Range("J104").Value = 31983
$ws.Range("L104").Value = 31983
$ws.Range("N104").Value = -38971

$ws.Range("H122").Value = 1076.6666
$ws.Range("I122").Value = 1562.5
$ws.Range("J122").Value = 833.75
$ws.Range("K122").Value = 4687.5
$ws.Range("L122").Value = 2501.25
$ws.Range("M122").Value = -2237.5
$ws.Range("N122").Value = -7401.25

$ws.Range("H126").Value = 2049.85
$ws.Range("I126").Value = 2131.158
$ws.Range("J126").Value = 505
$ws.Range("K126").Value = 6393.474
$ws.Range("L126").Value = 1515
$ws.Range("M126").Value = -3923.474
$ws.Range("N126").Value = -6455

$ws.Range("H132").Value = 11937.066
$ws.Range("I132").Value = 17949.777
$ws.Range("J132").Value = 2918
$ws.Range("K132").Value = 53849.33099999999
$ws.Range("L132").Value = 8754
$ws.Range("M132").Value = -51319.33099999999
$ws.Range("N132").Value = -13814

$ws.Range("H136").Value = 13862.4375
$ws.Range("I136").Value = 29872
$ws.Range("J136").Value = 1410.5555
$ws.Range("K136").Value = 89616
$ws.Range("L136").Value = 4231.666499999999
$ws.Range("M136").Value = -87066
$ws.Range("N136").Value = -9331.666499999999
